$d = $word.ActiveDocument

# The original document has a single paragraph:
#   "Hi, here is some text!!!"
# We need to split this into three paragraphs:
#   1. "Hi, here is some text!!!"
#   2. "Here is a title" (styled as Heading 1)
#   3. "And here is the rest of this text!!!"
# The trailing bookmark (_GoBack) should end up on the last paragraph.

$para = $d.Paragraphs(1)
$r = $para.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Here is a title")

$para2 = $d.Paragraphs(2)
$r2 = $para2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.Collapse(0)
$r2.InsertAfter("And here is the rest of this text!!!")

$d.Paragraphs(2).Style = "Heading 1"
